# Add "Mid Paper 1" / "Mid Paper 2" columns (F, G) to both mark sheets,
# populate the mid-term marks (mirroring the final Paper 1 / Paper 2 marks
# for the two students that already have scores), update each sheet's
# selection to the new columns, and make "Senior Five" the active tab.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Senior Six", "Senior Five")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F1").Value = "Mid Paper 1"
    $ws.Range("G1").Value = "Mid Paper 2"

    # Row 2 has no marks yet (same as the existing D2/E2 cells) -- touch the
    # cells so they materialise in the sheet without giving them a value.
    $ws.Range("F2").NumberFormat = "General"
    $ws.Range("G2").NumberFormat = "General"

    $ws.Range("F5").Value = 62
    $ws.Range("G5").Value = 54

    $ws.Range("F6").Value = 89
    $ws.Range("G6").Value = 78

    $null = $ws.Range("F1:G6").Select()
}

# "Senior Five" becomes the active sheet/tab.
$wb.Worksheets.Item("Senior Five").Activate()
